$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '323.41'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-2.01%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.31'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.76%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.725'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '8.87%'

# Row 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.25%'

# Row 6
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.507'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.35%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.610'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.15%'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.968'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.84%'

# Row 9
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.947'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.34%'

# Row 10
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9262'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.03%'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1246'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-5.86%'

# Row 12
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1952'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.92%'

# Row 13
$ws.Range("B13").Value = 'MCDex'
$ws.Range("C13").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.705'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '24.89%'

# Row 14
$ws.Range("B14").Value = 'MandalaExchangeToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09222'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.79%'

# Row 15
$ws.Range("B15").Value = 'BitrueCoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.03648'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.24%'

# Row 16
$ws.Range("B16").Value = 'BitMartToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.1050'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '9.66%'

# Row 17
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001286'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-3.30%'

# Row 18
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006180'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '4.64%'

# Row 19
$ws.Range("B19").Value = 'LEO'
$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.350'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.49%'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3538'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.28%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1372'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.34%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2451'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-4.34%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04414'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.16%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001263'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.55%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004571'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '6.18%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001152'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-3.24%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02510'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.58%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05331'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.88%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007447'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.20%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009573'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.19%'

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.68%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002120'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.92%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01071'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '4.78%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006779'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '1.72%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.16%'

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002294'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-7.60%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.16%'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.16%'
